$d = $word.ActiveDocument

# Locate the paragraph that reads "HomePage design - FeaturedPosts section"
# (the run text as seen by Word does not include the literal en-dash glyph
# distinctions, so a simple substring search on "FeaturedPosts" is reliable).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*FeaturedPosts*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'HomePage design - FeaturedPosts section' paragraph"
}

# --- Step 1: wrap "FeaturedPosts" in spell-check proofErr markers -----------
# Rebuild the paragraph's own content (minus its end-of-paragraph mark) so the
# "FeaturedPosts" run gets <w:proofErr w:type="spellStart"/> / spellEnd
# around it, matching the rest of the document's existing convention (the
# "HomePage" run already carries this marking and is left untouched).
$fullRange = $target.Range
$contentRange = $d.Range($fullRange.Start, $fullRange.End - 1)

$frag243 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Khmer UI" w:hAnsi="Khmer UI" w:cs="Khmer UI"/><w:lang w:val="en-US"/></w:rPr><w:t>HomePage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Khmer UI" w:hAnsi="Khmer UI" w:cs="Khmer UI"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> design – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Khmer UI" w:hAnsi="Khmer UI" w:cs="Khmer UI"/><w:lang w:val="en-US"/></w:rPr><w:t>FeaturedPosts</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Khmer UI" w:hAnsi="Khmer UI" w:cs="Khmer UI"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> section</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$contentRange.InsertXML($frag243)

# --- Step 2: insert a new paragraph right after it --------------------------
# "HomePage design - PostList section" -- same Khmer UI / en-US formatting,
# with "HomePage" proof-marked (as typed) and "PostList" left alone (it is
# not flagged as a spelling error by the author's checker).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*FeaturedPosts*") {
        $target = $p
        break
    }
}
$insertPoint = $d.Range($target.Range.End - 1, $target.Range.End - 1)

$frag244 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Khmer UI" w:hAnsi="Khmer UI" w:cs="Khmer UI"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Khmer UI" w:hAnsi="Khmer UI" w:cs="Khmer UI"/><w:lang w:val="en-US"/></w:rPr><w:t>HomePage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Khmer UI" w:hAnsi="Khmer UI" w:cs="Khmer UI"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> design – </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Khmer UI" w:hAnsi="Khmer UI" w:cs="Khmer UI"/><w:lang w:val="en-US"/></w:rPr><w:t>PostList</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Khmer UI" w:hAnsi="Khmer UI" w:cs="Khmer UI"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> section</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertPoint.InsertXML($frag244)

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
